# Update countries & provincias Spain
# Applies the COVID-19 dashboard data refresh:
#  - Catar's case count overtook Rumania and Emiratos Arabes Unidos, so it
#    moves up three rows (37 -> 35); Rumania and Emiratos each shift down one
#    row and keep their (now-stale) previous figures until their own refresh.
#  - A handful of other country rows get refreshed totals.
#  - The "datos actualizados" timestamp cell advances from 13:52 to 14:22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 14:22"

# --- Row 8: Alemania ---
$ws.Cells.Item(8, 4).Value = 117400
$ws.Cells.Item(8, 5).Value = 35232

# --- Row 17: Paises Bajos ---
$ws.Cells.Item(17, 2).Value = 38416
$ws.Cells.Item(17, 3).Value = 171
$ws.Cells.Item(17, 5).Value = 33600
$ws.Cells.Item(17, 7).Value = 48
$ws.Cells.Item(17, 8).Value = 4566

# --- Row 21: Portugal ---
$ws.Cells.Item(21, 2).Value = 24322
$ws.Cells.Item(21, 3).Value = 295
$ws.Cells.Item(21, 4).Value = 1389
$ws.Cells.Item(21, 5).Value = 21985
$ws.Cells.Item(21, 6).Value = 172
$ws.Cells.Item(21, 7).Value = 20
$ws.Cells.Item(21, 8).Value = 948

# --- Row 24: Suecia ---
$ws.Cells.Item(24, 2).Value = 19621
$ws.Cells.Item(24, 3).Value = 695
$ws.Cells.Item(24, 5).Value = 16261
$ws.Cells.Item(24, 6).Value = 524
$ws.Cells.Item(24, 7).Value = 81
$ws.Cells.Item(24, 8).Value = 2355

# --- Rows 35-37: Catar overtakes Rumania and Emiratos Arabes Unidos ---
# Row 35 becomes Catar, with its refreshed figures.
$ws.Cells.Item(35, 1).Value = "Catar"
$ws.Cells.Item(35, 2).Value = 11921
$ws.Cells.Item(35, 3).Value = 677
$ws.Cells.Item(35, 4).Value = 1134
$ws.Cells.Item(35, 5).Value = 10777
$ws.Cells.Item(35, 6).Value = 72
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = 10

# Row 36 becomes Rumania (shifted down from row 35, keeping its old figures).
$ws.Cells.Item(36, 1).Value = "Rumania"
$ws.Cells.Item(36, 2).Value = 11616
$ws.Cells.Item(36, 3).Value = 277
$ws.Cells.Item(36, 4).Value = 3404
$ws.Cells.Item(36, 5).Value = 7562
$ws.Cells.Item(36, 6).Value = 243
$ws.Cells.Item(36, 7).Value = 9
$ws.Cells.Item(36, 8).Value = 650

# Row 37 becomes Emiratos Arabes Unidos (shifted down from row 36).
$ws.Cells.Item(37, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(37, 2).Value = 11380
$ws.Cells.Item(37, 3).Value = 541
$ws.Cells.Item(37, 4).Value = 2181
$ws.Cells.Item(37, 5).Value = 9110
$ws.Cells.Item(37, 6).Value = 1
$ws.Cells.Item(37, 7).Value = 7
$ws.Cells.Item(37, 8).Value = 89

# --- Row 41: Dinamarca ---
$ws.Cells.Item(41, 4).Value = 6121
$ws.Cells.Item(41, 5).Value = 2296
$ws.Cells.Item(41, 6).Value = 66
$ws.Cells.Item(41, 7).Value = 7
$ws.Cells.Item(41, 8).Value = 434

# --- Row 61: Kazajistan ---
$ws.Cells.Item(61, 2).Value = 3019
$ws.Cells.Item(61, 3).Value = 184
$ws.Cells.Item(61, 4).Value = 754
$ws.Cells.Item(61, 5).Value = 2240

# --- Row 67: Croacia ---
$ws.Cells.Item(67, 2).Value = 2047
$ws.Cells.Item(67, 3).Value = 8
$ws.Cells.Item(67, 4).Value = 1232
$ws.Cells.Item(67, 5).Value = 752
$ws.Cells.Item(67, 7).Value = 4
$ws.Cells.Item(67, 8).Value = 63

# --- Row 105: Sri Lanka ---
$ws.Cells.Item(105, 2).Value = 596
$ws.Cells.Item(105, 3).Value = 8
$ws.Cells.Item(105, 5).Value = 455

# --- Row 192: San Cristobal y Nieves ---
$ws.Cells.Item(192, 4).Value = 4
$ws.Cells.Item(192, 5).Value = 11
